# edit.ps1
# Applies the changes described by the commit "Added brief descriptions of
# the companies, per advice." to the resume document.
#
# Summary of changes:
#  1. Summary paragraph: trim the first sentence and split the paragraph in
#     two, adding a new BodyText paragraph with extra bio detail.
#  2. Add a BlockText paragraph describing McGraw Hill right after the
#     "McGraw Hill ... Remote (since COVID), previously Seattle, WA" line.
#  3. Add a new bullet (numId 1002) under "Principal Cloud and Platform
#     Engineer (June 2020-January 2024)" describing the COVID-19 pandemic
#     response work.
#  4. Add a BlockText paragraph describing WePay right after the
#     "WePay ... Redwood City, CA" line.
#  5. Add a BlockText paragraph describing Amazon Web Services right after
#     the "Amazon Web Services ... Seattle, WA" line.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: run Find.Execute and fail loudly if nothing was found, so that
# a broken anchor doesn't silently turn into a no-op.
# ---------------------------------------------------------------------
function Find-RequiredText($range, [string]$searchText, [string]$replaceText = "", [int]$replaceMode = 0) {
    $found = $range.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, $replaceMode)
    if (-not $found) {
        throw "Find-RequiredText: could not find required text: $searchText"
    }
}

# ---------------------------------------------------------------------
# Helper: insert a literal OOXML paragraph fragment immediately after the
# paragraph whose end is at the (already collapsed, end-of-range) position
# $insertPos. Works by creating an empty paragraph via InsertParagraphAfter,
# locating that freshly created (now-empty) paragraph inside the
# Paragraphs collection by its start offset, and overwriting it with the
# caller-supplied <w:p>...</w:p> XML via Range.InsertXML.
# ---------------------------------------------------------------------
function Insert-ParagraphXmlAfter($range, [string]$paragraphXml) {
    $range.Collapse(0) | Out-Null
    $insertPos = $range.End
    $range.InsertParagraphAfter()

    $paras = $d.Paragraphs
    $target = $null
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        if ($p.Range.Start -eq ($insertPos + 1)) {
            $target = $p.Range
            break
        }
    }
    if ($null -eq $target) {
        throw "Insert-ParagraphXmlAfter: could not locate newly inserted empty paragraph at position $($insertPos + 1)"
    }

    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
        $paragraphXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $target.InsertXML($xml) | Out-Null
}

# ---------------------------------------------------------------------
# 1. Summary paragraph: shorten, then add a second BodyText paragraph.
# ---------------------------------------------------------------------
$r = $d.Content
Find-RequiredText $r `
    "Ryan Parman is a cloud-native engineering leader with over 25 years of experience, who specializes in technical leadership, software development, site reliability engineering, and cybersecurity for the modern web. Excels at listening, adapting, and driving continuous improvement. Delivers exceptional work, builds impactful solutions, and elevates team performance." `
    "Ryan Parman is a cloud-native engineering leader, who specializes in technical leadership, software development, site reliability engineering, and cybersecurity for the modern web. Excels at listening, adapting, and driving continuous improvement." `
    2

$bodyText1 = '<w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">Small business owner, two-time startup founder, and creator of two open-source projects with millions of users each. Ryan has a proven track record of building high-quality software, delivering impactful solutions, and elevating team performance.</w:t></w:r></w:p>'
Insert-ParagraphXmlAfter $r $bodyText1

# ---------------------------------------------------------------------
# 2. McGraw Hill description (BlockText), after the "... Seattle, WA" line
#    for the McGraw Hill entry.
# ---------------------------------------------------------------------
$r = $d.Content
Find-RequiredText $r "— Remote (since COVID), previously Seattle, WA"

$mcgrawHill = '<w:p><w:pPr><w:pStyle w:val="BlockText"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">McGraw Hill is a</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve">learning science</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">company which produces textbooks, digital learning tools, and adaptive technology to enhance learning. It is one of the &#8220;big three&#8221; educational publishers in the U.S, and was acquired by Platinum Equity 2021.</w:t></w:r>' +
    '</w:p>'
Insert-ParagraphXmlAfter $r $mcgrawHill

# ---------------------------------------------------------------------
# 3. New bullet under "Principal Cloud and Platform Engineer (June
#    2020-January 2024)" about the COVID-19 transition to online learning.
# ---------------------------------------------------------------------
$r = $d.Content
Find-RequiredText $r "Principal Cloud and Platform Engineer (June 2020—January 2024)"

$covidBullet = '<w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1002"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">As every school in America transitioned to online learning during the COVID-19 lockdowns, I was the technical/development lead on the team who supported all SRE and product engineering teams, working on core platforms and services.</w:t></w:r></w:p>'
Insert-ParagraphXmlAfter $r $covidBullet

# ---------------------------------------------------------------------
# 4. WePay description (BlockText), after the "... Redwood City, CA" line.
# ---------------------------------------------------------------------
$r = $d.Content
Find-RequiredText $r "— Redwood City, CA"

$wepay = '<w:p><w:pPr><w:pStyle w:val="BlockText"/></w:pPr><w:r><w:t xml:space="preserve">WePay is an online payment service provider which provides &#8220;payments for platforms&#8221;, where examples of platforms are GoFundMe, Care.com, and Xbox. It was acquired by JPMorgan Chase in October 2017.</w:t></w:r></w:p>'
Insert-ParagraphXmlAfter $r $wepay

# ---------------------------------------------------------------------
# 5. Amazon Web Services description (BlockText), after the "... Seattle,
#    WA" line for the AWS entry.
# ---------------------------------------------------------------------
$r = $d.Content
Find-RequiredText $r "— Seattle, WA"

$aws = '<w:p><w:pPr><w:pStyle w:val="BlockText"/></w:pPr><w:r><w:t xml:space="preserve">Amazon Web Services provides on-demand cloud computing platforms and APIs to individuals, companies, and governments, on a metered, pay-as-you-go basis.</w:t></w:r></w:p>'
Insert-ParagraphXmlAfter $r $aws

Write-Host "All edits applied."
